# edit.ps1
#
# Implements commit "feat: add 2022-Q3 data":
#   - Inserts a new worksheet "2022-Q3" right after "总计" (i.e. before the
#     existing "2022-Q2" sheet), populated with the fund-holdings table for
#     that quarter.
#   - Updates the "总计" (summary) sheet: adds a new leading row for
#     "2022-Q3" (holdings count 33, market value 10.38) and shifts every
#     previously-existing quarter row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a cell, forcing text-typed values (kind "s") to stay text
# even when they look numeric (fund codes like "013082", percentages like
# "89.00") by using Excel's leading-apostrophe text prefix. Numeric-typed
# values (kind "n") are written as real numbers.
# ---------------------------------------------------------------------------
function Set-CellTyped {
    param($ws, $addr, $kind, $val)
    if ($kind -eq "n") {
        $ws.Range($addr).Value = [double]$val
    } else {
        $ws.Range($addr).Value = "'" + $val
    }
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet immediately before "2022-Q2"
#    (i.e. right after "总计"), matching the target sheet order:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# Grab the header/index-column formatting (bold, bordered, centered) from the
# existing "2022-Q2" sheet before we touch anything, so the new sheet's
# header row (B1:H1) and index column (A2:A34) can reuse the very same
# cell style instead of Excel fabricating a brand-new one.
$q2Sheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2:A34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fund-holdings data for 2022-Q3 (34 rows incl. header, columns A-H).
#    Each row entry is: rowNumber, then (column, kind, value) tuples, where
#    kind "n" = numeric, "s" = text.
# ---------------------------------------------------------------------------
$q3rows = @(
    @(1, @("B","s","基金代码"),@("C","s","基金名称"),@("D","s","基金规模"),@("E","s","股票总仓位"),@("F","s","仓位占比"),@("G","s","持有市值(亿元)"),@("H","s","仓位排名")),
    @(2, @("A","n","0"),@("B","s","516970"),@("C","s","广发中证基建工程ETF"),@("D","s","89.00"),@("E","s","99.56"),@("F","s","6.03"),@("G","s","5.3667"),@("H","n","4")),
    @(3, @("A","n","1"),@("B","s","515900"),@("C","s","博时中证央企创新驱动ETF"),@("D","s","36.57"),@("E","s","98.62"),@("F","s","2.94"),@("G","s","1.0752"),@("H","n","7")),
    @(4, @("A","n","2"),@("B","s","165525"),@("C","s","信诚中证基建工程指数（LOF）"),@("D","s","10.82"),@("E","s","94.26"),@("F","s","5.71"),@("G","s","0.6178"),@("H","n","4")),
    @(5, @("A","n","3"),@("B","s","002670"),@("C","s","万家沪深300指数增强A"),@("D","s","20.85"),@("E","s","94.06"),@("F","s","2.48"),@("G","s","0.5171"),@("H","n","3")),
    @(6, @("A","n","4"),@("B","s","516950"),@("C","s","银华中证基建ETF"),@("D","s","11.07"),@("E","s","97.93"),@("F","s","4.34"),@("G","s","0.4804"),@("H","n","5")),
    @(7, @("A","n","5"),@("B","s","515600"),@("C","s","广发中证央企创新驱动ETF"),@("D","s","14.84"),@("E","s","98.78"),@("F","s","2.94"),@("G","s","0.4363"),@("H","n","7")),
    @(8, @("A","n","6"),@("B","s","515680"),@("C","s","嘉实中证央企创新驱动ETF"),@("D","s","14.64"),@("E","s","99.23"),@("F","s","2.93"),@("G","s","0.4290"),@("H","n","7")),
    @(9, @("A","n","7"),@("B","s","002671"),@("C","s","万家沪深300指数增强C"),@("D","s","10.38"),@("E","s","94.06"),@("F","s","2.48"),@("G","s","0.2574"),@("H","n","3")),
    @(10, @("A","n","8"),@("B","s","013082"),@("C","s","信诚中证基建工程指数（LOF）C"),@("D","s","3.48"),@("E","s","94.26"),@("F","s","5.71"),@("G","s","0.1987"),@("H","n","4")),
    @(11, @("A","n","9"),@("B","s","159635"),@("C","s","华夏中证基建ETF"),@("D","s","3.40"),@("E","s","99.03"),@("F","s","4.39"),@("G","s","0.1493"),@("H","n","5")),
    @(12, @("A","n","10"),@("B","s","159619"),@("C","s","国泰中证基建ETF"),@("D","s","3.30"),@("E","s","98.76"),@("F","s","4.38"),@("G","s","0.1445"),@("H","n","5")),
    @(13, @("A","n","11"),@("B","s","160135"),@("C","s","南方中证高铁产业指数（LOF）"),@("D","s","1.84"),@("E","s","95.01"),@("F","s","7.85"),@("G","s","0.1444"),@("H","n","5")),
    @(14, @("A","n","12"),@("B","s","159974"),@("C","s","富国中证央企创新驱动ETF"),@("D","s","4.89"),@("E","s","99.47"),@("F","s","2.95"),@("G","s","0.1443"),@("H","n","7")),
    @(15, @("A","n","13"),@("B","s","004497"),@("C","s","前海开源多元策略灵活配置混合C"),@("D","s","1.68"),@("E","s","93.04"),@("F","s","4.83"),@("G","s","0.0811"),@("H","n","5")),
    @(16, @("A","n","14"),@("B","s","169108"),@("C","s","东方红均衡优选两年定期开放混合"),@("D","s","7.09"),@("E","s","25.64"),@("F","s","1.04"),@("G","s","0.0737"),@("H","n","5")),
    @(17, @("A","n","15"),@("B","s","160639"),@("C","s","鹏华中证高铁产业指数（LOF）A"),@("D","s","0.75"),@("E","s","94.62"),@("F","s","7.79"),@("G","s","0.0584"),@("H","n","5")),
    @(18, @("A","n","16"),@("B","s","004496"),@("C","s","前海开源多元策略灵活配置混合A"),@("D","s","0.91"),@("E","s","93.04"),@("F","s","4.83"),@("G","s","0.0440"),@("H","n","5")),
    @(19, @("A","n","17"),@("B","s","011471"),@("C","s","鹏华致远成长混合A"),@("D","s","1.84"),@("E","s","65.59"),@("F","s","1.92"),@("G","s","0.0353"),@("H","n","9")),
    @(20, @("A","n","18"),@("B","s","517090"),@("C","s","国泰富时中国国企开放共赢ETF"),@("D","s","0.60"),@("E","s","91.75"),@("F","s","3.68"),@("G","s","0.0221"),@("H","n","8")),
    @(21, @("A","n","19"),@("B","s","011050"),@("C","s","天弘裕新混合A"),@("D","s","1.61"),@("E","s","21.83"),@("F","s","1.34"),@("G","s","0.0216"),@("H","n","9")),
    @(22, @("A","n","20"),@("B","s","000423"),@("C","s","前海开源事件驱动混合A"),@("D","s","0.45"),@("E","s","89.39"),@("F","s","4.45"),@("G","s","0.0200"),@("H","n","10")),
    @(23, @("A","n","21"),@("B","s","011048"),@("C","s","天弘恒新混合A"),@("D","s","0.91"),@("E","s","25.66"),@("F","s","2.06"),@("G","s","0.0187"),@("H","n","5")),
    @(24, @("A","n","22"),@("B","s","517180"),@("C","s","南方富时中国国企开放共赢ETF"),@("D","s","0.36"),@("E","s","97.73"),@("F","s","2.33"),@("G","s","0.0084"),@("H","n","10")),
    @(25, @("A","n","23"),@("B","s","011051"),@("C","s","天弘裕新混合C"),@("D","s","0.59"),@("E","s","21.83"),@("F","s","1.34"),@("G","s","0.0079"),@("H","n","9")),
    @(26, @("A","n","24"),@("B","s","159719"),@("C","s","平安富时中国国企开放共赢ETF"),@("D","s","0.34"),@("E","s","94.14"),@("F","s","2.29"),@("G","s","0.0078"),@("H","n","10")),
    @(27, @("A","n","25"),@("B","s","015678"),@("C","s","鹏华中证高铁产业指数（LOF）C"),@("D","s","0.06"),@("E","s","94.62"),@("F","s","7.79"),@("G","s","0.0047"),@("H","n","5")),
    @(28, @("A","n","26"),@("B","s","001185"),@("C","s","安信动态策略灵活配置混合A"),@("D","s","0.27"),@("E","s","32.82"),@("F","s","1.40"),@("G","s","0.0038"),@("H","n","9")),
    @(29, @("A","n","27"),@("B","s","002029"),@("C","s","安信动态策略灵活配置混合C"),@("D","s","0.24"),@("E","s","32.82"),@("F","s","1.40"),@("G","s","0.0034"),@("H","n","9")),
    @(30, @("A","n","28"),@("B","s","011049"),@("C","s","天弘恒新混合C"),@("D","s","0.16"),@("E","s","25.66"),@("F","s","2.06"),@("G","s","0.0033"),@("H","n","5")),
    @(31, @("A","n","29"),@("B","s","004360"),@("C","s","创金合信量化核心混合C"),@("D","s","0.21"),@("E","s","91.03"),@("F","s","1.33"),@("G","s","0.0028"),@("H","n","9")),
    @(32, @("A","n","30"),@("B","s","001865"),@("C","s","前海开源事件驱动混合C"),@("D","s","0.05"),@("E","s","89.39"),@("F","s","4.45"),@("G","s","0.0022"),@("H","n","10")),
    @(33, @("A","n","31"),@("B","s","004359"),@("C","s","创金合信量化核心混合A"),@("D","s","0.16"),@("E","s","91.03"),@("F","s","1.33"),@("G","s","0.0021"),@("H","n","9")),
    @(34, @("A","n","32"),@("B","s","011472"),@("C","s","鹏华致远成长混合C"),@("D","s","0.07"),@("E","s","65.59"),@("F","s","1.92"),@("G","s","0.0013"),@("H","n","9"))
)

foreach ($rowdef in $q3rows) {
    $rnum = $rowdef[0]
    for ($i = 1; $i -lt $rowdef.Length; $i++) {
        $cell = $rowdef[$i]
        $col = $cell[0]
        $kind = $cell[1]
        $val = $cell[2]
        Set-CellTyped $newSheet "$col$rnum" $kind $val
    }
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert the new 2022-Q3 row at the top
#    of the data (row 2) and shift the rest down by one row. Rewriting all
#    data rows outright (rather than relying on Rows.Insert, which pulls in
#    inconsistent formatting) keeps the existing header/index-column style
#    intact for every row.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

$totalRows = @(
    @(0, "2022-Q3", 33, 10.38),
    @(1, "2022-Q2", 28, 11.47),
    @(2, "2022-Q1", 21, 8.220000000000001),
    @(3, "2021-Q4", 16, 4.48),
    @(4, "2021-Q3", 19, 8.67),
    @(5, "2021-Q2", 16, 7.76),
    @(6, "2021-Q1", 10, 4.7),
    @(7, "2020-Q4", 7, 4.78)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $totalWs.Cells.Item($r, 1).Value = [double]$row[0]
    $totalWs.Cells.Item($r, 2).Value = "'" + $row[1]
    $totalWs.Cells.Item($r, 3).Value = [double]$row[2]
    $totalWs.Cells.Item($r, 4).Value = [double]$row[3]
}

# Row 9 is brand new (the table grew from 8 to 9 data rows); give its index
# cell (A9) the same style as the other index-column cells by copying it
# from A2 (which already carries that formatting).
$totalWs.Range("A2").Copy()
$totalWs.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
